$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add - [AI / Item, Equipment] MarketAI ItemSlot additions: BIGBOX_COUNT / SMALLBOX_COUNT
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "BIGBOX_COUNT"
$ws.Range("C21").Value = 100

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "SMALLBOX_COUNT"
$ws.Range("C22").Value = 25

# Match the vertically-centered style used by the rest of column A/B
$ws.Range("A21:B21").VerticalAlignment = -4108
$ws.Range("A22").VerticalAlignment = -4108

# Leave the new rows selected, as in the authored change
$ws.Range("A21:C22").Select()
